$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Record Number"
$ws.Range("B1").Value = "Latitude"
$ws.Range("C1").Value = "Longitude"
$ws.Range("D1").Value = "Summary"

# Row 2 - first accident data point
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 39.2
$ws.Range("C2").Value = -86.4
$ws.Range("D2").Value = "This is a summary of the first accident data point. "

# Row 3 - second accident data point
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = -87
$ws.Range("D3").Value = "This is a summary of the second accident data point. "

# Select D3 to match the saved selection state
$ws.Range("D3").Select()

# Set column A width to (approximately) match the bestFit width seen in the
# target file for the "Record Number" header text.
$ws.Columns.Item(1).ColumnWidth = 13
